$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# October 6th (row 7) daily log entries.
# Copy the existing conditional-style formatting (Good/Neutral/Bad cell
# styles, applied per-column in the rows above) onto row 7 before writing
# the new values, so the new cells pick up the same style indices already
# used elsewhere in the sheet instead of creating new ones.

$ws.Range("B6").Copy() | Out-Null
$ws.Range("B7").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("C6").Copy() | Out-Null
$ws.Range("C7").PasteSpecial(-4122) | Out-Null

$ws.Range("D6").Copy() | Out-Null
$ws.Range("D7").PasteSpecial(-4122) | Out-Null

$ws.Range("E6").Copy() | Out-Null
$ws.Range("E7").PasteSpecial(-4122) | Out-Null

$ws.Range("F6").Copy() | Out-Null
$ws.Range("F7").PasteSpecial(-4122) | Out-Null

$ws.Range("G6").Copy() | Out-Null
$ws.Range("G7").PasteSpecial(-4122) | Out-Null

$ws.Range("H6").Copy() | Out-Null
$ws.Range("H7").PasteSpecial(-4122) | Out-Null

$ws.Range("I6").Copy() | Out-Null
$ws.Range("I7").PasteSpecial(-4122) | Out-Null

# Column J is "Spécialisations" - on row 7 it gets the Neutral style
# (same style as column F), unlike rows 2-6 where it uses the Good style.
$ws.Range("F6").Copy() | Out-Null
$ws.Range("J7").PasteSpecial(-4122) | Out-Null

$ws.Range("K6").Copy() | Out-Null
$ws.Range("K7").PasteSpecial(-4122) | Out-Null

$ws.Range("L6").Copy() | Out-Null
$ws.Range("L7").PasteSpecial(-4122) | Out-Null

$ws.Range("M6").Copy() | Out-Null
$ws.Range("M7").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Now fill in the actual values for October 6th, in the same order they
# were originally typed (this determines the order new entries land in
# the shared-string table).
$ws.Range("J7").Value = "Archer"
$ws.Range("B7").Value = "Guild Artisan"
$ws.Range("C7").Value = "Force"
$ws.Range("F7").Value = "Dans la rue"
$ws.Range("G7").Value = "Bricolage magique"
$ws.Range("H7").Value = "Emboîter le pas"
$ws.Range("I7").Value = "Attaque neutralisante"
$ws.Range("E7").Value = "Borgne, Chance extraordinaire"
